$d = $word.ActiveDocument

# The last paragraph in the "KNOWN BUGS" list currently holds only the
# _GoBack bookmark. Turn it into a new bulleted bug entry:
#   "OSC control can " + "show" + <bookmark> + " bugs when used with multiple control devices"
$lastIndex = $d.Paragraphs.Count
$bugsPara = $d.Paragraphs.Item($lastIndex)
$refPara = $d.Paragraphs.Item($lastIndex - 1)

$bugsRange = $bugsPara.Range

# Insert the trailing text after the bookmark, then the two leading runs
# before it, so the bookmark stays put between "show" and " bugs ...".
$bugsRange.InsertAfter(" bugs when used with multiple control devices")
$bugsRange.InsertBefore("show")
$bugsRange.InsertBefore("OSC control can ")

# Make it a list item matching the other "KNOWN BUGS" entries (style
# "Paragraphedeliste" / numId 3), reusing the existing list instance
# instead of starting a new one.
$bugsPara.Style = "List Paragraph"
$bugsPara.Range.ListFormat.ApplyListTemplateWithLevel($refPara.Range.ListFormat.ListTemplate, $true, 0, $false, 0)

# Append a new, empty trailing paragraph after it.
$bugsPara2 = $d.Paragraphs.Item($lastIndex)
$endRange = $d.Range($bugsPara2.Range.End, $bugsPara2.Range.End)
$endRange.Text = "`r"

# Strip the formatting the new paragraph mark inherited, so it is a bare
# empty paragraph (no list/style carried over).
$newPara = $d.Paragraphs.Item($lastIndex + 1)
$newPara.Range.ListFormat.RemoveNumbers()
$newPara.Style = "Normal"
